$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 207, shifting existing rows 207+ down by one.
$ws.Rows.Item(207).Insert()

# Populate the new row 207 with the data from the commit.
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").NumberFormat = $ws.Range("D206").NumberFormat
$ws.Range("D207").Value = 44476
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100108
$ws.Range("H207").Value = "Tropicales y subtropicales"
$ws.Range("I207").Value = 100108006
$ws.Range("J207").Value = "Plátano"
$ws.Range("K207").Value = "Sin especificar"
$ws.Range("L207").Value = "Primera Pintón"
$ws.Range("M207").Value = 1000
$ws.Range("N207").Value = 22000
$ws.Range("O207").Value = 23000
$ws.Range("P207").Value = 22500
$ws.Range("Q207").Value = "$/caja 20 kilos"
$ws.Range("R207").Value = "Ecuador"
$ws.Range("S207").Value = 1125
$ws.Range("T207").Value = 20
